$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) block labels
$ws.Range("A1").Value = "kitchens_1"
$ws.Range("B1").Value = "living_rooms_1"
$ws.Range("D1").Value = "kitchens_2"
$ws.Range("E1").Value = "bedrooms_2"

# Update block-order indicator matrix (rows 2-7)
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 1

$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 1

$ws.Range("A4").Value = 0
$ws.Range("B4").Value = 1

$ws.Range("A6").Value = 1
$ws.Range("E6").Value = 0
